# Update the answer values in the first (and only) table on the page.
# Each cell is addressed directly by (row, column) so the handful of
# duplicate "old" values elsewhere in the sheet aren't affected.
$d = $word.ActiveDocument
$t = $d.Tables(1)
$t.Cell(1,1).Range.Text = "53÷9=5, 8"
$t.Cell(1,2).Range.Text = "54÷7=7, 5"
$t.Cell(1,3).Range.Text = "62÷4=15, 2"
$t.Cell(1,4).Range.Text = "63÷6=10, 3"
$t.Cell(1,5).Range.Text = "17÷6=2, 5"
$t.Cell(5,1).Range.Text = "12÷8=1, 4"
$t.Cell(5,2).Range.Text = "41÷2=20, 1"
$t.Cell(5,3).Range.Text = "29÷8=3, 5"
$t.Cell(5,4).Range.Text = "64÷8=8, 0"
$t.Cell(5,5).Range.Text = "96÷3=32, 0"
$t.Cell(9,1).Range.Text = "72÷7=10, 2"
$t.Cell(9,2).Range.Text = "90÷7=12, 6"
$t.Cell(9,3).Range.Text = "69÷4=17, 1"
$t.Cell(9,4).Range.Text = "91÷4=22, 3"
$t.Cell(9,5).Range.Text = "20÷5=4, 0"
$t.Cell(13,1).Range.Text = "22÷2=11, 0"
$t.Cell(13,3).Range.Text = "35÷2=17, 1"
$t.Cell(13,4).Range.Text = "27÷2=13, 1"
$t.Cell(13,5).Range.Text = "18÷6=3, 0"
$t.Cell(17,1).Range.Text = "54÷3=18, 0"
$t.Cell(17,2).Range.Text = "64÷4=16, 0"
$t.Cell(17,3).Range.Text = "58÷2=29, 0"
$t.Cell(17,4).Range.Text = "25÷7=3, 4"
$t.Cell(17,5).Range.Text = "20÷8=2, 4"
Write-Output "done"